$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 92, shifting existing rows 92-116 down to 93-117.
$ws.Rows.Item(92).Insert()

# Populate the new row 92 with the new weekly price record.
$ws.Range("A92").Value = 4
$ws.Range("B92").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C92").Value = "Los Lagos"
$ws.Range("D92").Value = (Get-Date -Year 2022 -Month 10 -Day 21 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E92").Value = 10
$ws.Range("F92").Value = "Fruta"
$ws.Range("G92").Value = 100104
$ws.Range("H92").Value = "Frutos de pepita"
$ws.Range("I92").Value = 100104003
$ws.Range("J92").Value = "Membrillo"
$ws.Range("K92").Value = "Champion"
$ws.Range("L92").Value = "Primera"
$ws.Range("M92").Value = 120
$ws.Range("N92").Value = 14000
$ws.Range("O92").Value = 15000
$ws.Range("P92").Value = 14500
$ws.Range("Q92").Value = "$/caja 18 kilos granel"
$ws.Range("R92").Value = "Región de O'Higgins"
$ws.Range("S92").Value = 806
$ws.Range("T92").Value = 18
